$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 299
$ws.Range("F3").Value = 496
$ws.Range("F4").Value = 42
$ws.Range("F5").Value = 48
$ws.Range("F7").Value = 1257
$ws.Range("F8").Value = 387
$ws.Range("F9").Value = 245
$ws.Range("F10").Value = 360
$ws.Range("F11").Value = 8262
$ws.Range("F13").Value = 10019
$ws.Range("F14").Value = 86
$ws.Range("F20").Value = 224
$ws.Range("F24").Value = 60
$ws.Range("F27").Value = 1713
$ws.Range("F28").Value = 47
$ws.Range("F29").Value = 498
$ws.Range("F30").Value = 320
$ws.Range("F31").Value = 277
$ws.Range("F33").Value = 555
$ws.Range("F34").Value = 1030
$ws.Range("F37").Value = 1404
$ws.Range("F38").Value = 422
$ws.Range("F39").Value = 328
$ws.Range("F42").Value = 489
$ws.Range("F43").Value = 316
$ws.Range("F44").Value = 67
$ws.Range("F46").Value = 111
$ws.Range("F48").Value = 38
$ws.Range("F49").Value = 43

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 15
$ws.Range("F5").Value = 100
$ws.Range("F12").Value = 10
$ws.Range("F15").Value = 57
$ws.Range("F19").Value = 35

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 208
$ws.Range("F3").Value = 2773
$ws.Range("F4").Value = 332
$ws.Range("F5").Value = 196

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 299
$ws.Range("F3").Value = 496
$ws.Range("F5").Value = 332
$ws.Range("F6").Value = 196
$ws.Range("F7").Value = 42
$ws.Range("F8").Value = 48
$ws.Range("F10").Value = 1257
$ws.Range("F11").Value = 387
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 245
$ws.Range("F15").Value = 100
$ws.Range("F16").Value = 8262
$ws.Range("F18").Value = 10019
$ws.Range("F24").Value = 1713
$ws.Range("F25").Value = 47
$ws.Range("F26").Value = 320
$ws.Range("F27").Value = 277
$ws.Range("F30").Value = 555
$ws.Range("F32").Value = 10
$ws.Range("F36").Value = 1404
$ws.Range("F37").Value = 422
$ws.Range("F38").Value = 57
$ws.Range("F39").Value = 328
$ws.Range("F41").Value = 489
$ws.Range("F42").Value = 316
$ws.Range("F43").Value = 67
$ws.Range("F45").Value = 35
$ws.Range("F48").Value = 38
$ws.Range("F49").Value = 43
